$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Hola"
$ws.Range("A3").Value = "H"

$ws.Range("A3").Select()
